$wb = $excel.ActiveWorkbook

# --- Sheet "Feuil1" edits ---
$ws1 = $wb.Worksheets.Item("Feuil1")

# Row 18: N18 becomes a static value (400.0) instead of formula; P18 16 -> 4; U18 gets new text
$ws1.Range("N18").Value = 400.0
$ws1.Range("P18").Value = 4.0
$ws1.Range("U18").Value = "متكفل بها"

# Row 19: N19 becomes static value (600.0) instead of formula; P19 7 -> 2; U19 gets new text
$ws1.Range("N19").Value = 600.0
$ws1.Range("P19").Value = 2.0
$ws1.Range("U19").Value = "متكفل بها"

# Row 20: N20, O20, P20 cleared
$ws1.Range("N20").ClearContents()
$ws1.Range("O20").ClearContents()
$ws1.Range("P20").ClearContents()

# Row 21: N21, O21, P21 cleared
$ws1.Range("N21").ClearContents()
$ws1.Range("O21").ClearContents()
$ws1.Range("P21").ClearContents()

# Row 23: N23 formula range narrows from N18:N22 to N18:N19
$ws1.Range("N23").Formula = "=SUM(N18:N19)"

# Row 25: N25 23300 -> 1000
$ws1.Range("N25").Value = 1000.0

# Row 38: L38 text changes ('2020/02/16' -> '2021/07/19'). L38 is a date-formatted
# cell (s=135, numFmtId 14) but the target content is a literal text string, not a
# real date value - so we must avoid Excel's "looks like a date" auto-conversion.
# Route the text through TEXT() + copy/paste-values (using the already-blank,
# same-styled M38 cell as scratch space) so the destination keeps its original
# style/number-format untouched while acquiring a plain text value.
$tmp38 = $ws1.Range("M38")
$tmp38.Formula = "=TEXT(44396,""yyyy/mm/dd"")"
$tmp38.Copy()
$ws1.Range("L38").PasteSpecial(-4163)  # xlPasteValues
$tmp38.ClearContents()

# --- Sheet "Feuil2" edits ---
$ws2 = $wb.Worksheets.Item("Feuil2")

# Row 9
$ws2.Range("B9").Value = "255"
$ws2.Range("D9").Value = 2.0
$ws2.Range("F9").Value = 4.0
$ws2.Range("M9").Value = "06:00"

# N9 ('20/02/2020' -> '06/11/2019') and P9 ('13/02/2020' -> '04/11/2019') are
# date-formatted cells (yyyy-dd-mm) holding literal text, not real dates - same
# quote-prefix problem as L38 above. Use the blank Q9 cell (General format) as
# scratch space to build the text via TEXT() then paste-values back in, keeping
# N9/P9's own style/number-format untouched.
$tmp9 = $ws2.Range("Q9")
$tmp9.Formula = "=TEXT(43775,""dd/mm/yyyy"")"
$tmp9.Copy()
$ws2.Range("N9").PasteSpecial(-4163)  # xlPasteValues
$tmp9.Formula = "=TEXT(43773,""dd/mm/yyyy"")"
$tmp9.Copy()
$ws2.Range("P9").PasteSpecial(-4163)  # xlPasteValues
$tmp9.ClearContents()

# Row 10: whole row content cleared
$ws2.Range("B10").Value = ""
$ws2.Range("C10").ClearContents()
$ws2.Range("E10").ClearContents()
$ws2.Range("L10").ClearContents()
$ws2.Range("M10").ClearContents()
$ws2.Range("N10").ClearContents()
$ws2.Range("O10").ClearContents()
$ws2.Range("P10").ClearContents()
$ws2.Range("R10").ClearContents()
$ws2.Range("S10").ClearContents()
$ws2.Range("T10").ClearContents()
